# Fruta / hortaliza, semanal
# Rotate the weekly snapshot data (Fecha/Variedad/Volumen/Precios) among the
# existing rows 2-14 (row 10 keeps its own data), per the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for each row: Fecha(D), Variedad(H), Volumen(J), Precio minimo(K),
# Precio maximo(L), Precio promedio ponderado(M), Precio $/Kg(P)
$rows = @{
    2  = @(44371, "Sin especificar", 80,  7000, 8000, 7375, 7375)
    3  = @(44309, "Sin especificar", 50,  8000, 9000, 8500, 8500)
    4  = @(44414, "Sin especificar", 100, 6000, 7000, 6500, 6500)
    5  = @(44539, "Americana (o)",   160, 6500, 7000, 6750, 6750)
    6  = @(44497, "Sin especificar", 160, 5000, 6000, 5500, 5500)
    7  = @(44263, "Sin especificar", 100, 7000, 8000, 7500, 7500)
    8  = @(44259, "Sin especificar", 80,  4000, 4500, 4250, 4250)
    9  = @(44410, "Sin especificar", 100, 5500, 6000, 5750, 5750)
    11 = @(44636, "Americana (o)",   60,  8000, 9000, 8500, 8500)
    12 = @(44575, "Sin especificar", 160, 6500, 7000, 6750, 6750)
    13 = @(44281, "Sin especificar", 100, 5000, 6000, 5500, 5500)
    14 = @(44559, "Americana (o)",   100, 5000, 6000, 5500, 5500)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 8).Value  = $vals[1]   # H - Variedad
    $ws.Cells.Item($r, 10).Value = $vals[2]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals[3]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals[4]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals[5]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals[6]   # P - Precio $/Kg
}
